$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-XRow($row, $aText, $fText, $fWrap) {
    $rangeA = "A$row`:E$row"
    $rangeF = "F$row`:I$row"
    $ws.Range($rangeA).HorizontalAlignment = -4108
    $ws.Range($rangeF).HorizontalAlignment = -4108
    if ($fWrap) {
        $ws.Range($rangeF).WrapText = $true
    }
    $ws.Range("A$row").Value = $aText
    $ws.Range("F$row").Value = $fText
    $mergeA = "A$row`:E" + ($row + 1)
    $mergeF = "F$row`:I" + ($row + 1)
    $ws.Range($mergeA).Merge()
    $ws.Range($mergeF).Merge()
}

Set-XRow 63 "profile image" '//img[@class="_3X2gOt"]' $false
Set-XRow 65 "form(in travel)" '//input[@class="_1w3ZZo _1YBGQV _2EjOJB lZd1T6 _2vegSu _2mFmU7"]' $true
Set-XRow 67 "depart on" '//input[@class="_1w3ZZo _2gKfhi _2mFmU7" and @name="0-datefrom"]' $true
Set-XRow 69 "return on" '//input[@class="_1w3ZZo _2gKfhi _2mFmU7" and @name="0-dateto"]' $true
Set-XRow 71 "search" '//button[@class="_2KpZ6l _1QYQF8 _3dESVI"]' $true
Set-XRow 73 "bakground image" '//div[@class="aCgX3e"]' $false

$ws.Range("I78").Select()
$excel.ActiveWindow.ScrollRow = 62
$excel.ActiveWindow.ScrollColumn = 1
